$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (Oct 10 2020 vs Kings XI Punjab) - shifts remaining rows up
$ws.Rows.Item(2).Delete()

# Delete what is now row 3 (originally row 4: Sep 30 2020 vs Rajasthan Royals)
$ws.Rows.Item(3).Delete()
